$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last row (row 53) so the table shrinks from 53 to 52 data+header rows,
# matching the new dimension A1:E52.
$ws.Rows.Item(53).Delete()

# Target data for rows 2..52 (columns A=date serial, B=y_0, C=y_0_forecast, D=y_1, E=y_1_forecast).
# Blank string means the cell should have no value (cleared).
$data = @"
2,39583,2008,,2009,
3,39765,2008,,2009,
4,39948,2009,,2010,
5,40130,2009,0.1715429114845124,2010,
6,40310,2010,,2011,
7,40494,2010,0.8004663283405655,2011,
8,40676,2011,,2012,
9,40862,2011,5.253783907501819,2012,
10,41044,2012,,2013,
11,41228,2012,3.522405026196918,2013,0.5784444854042281
12,41409,2013,1.133560223479058,2014,1.985690391709771
13,41592,2013,1.656063945467268,2014,2.529895848567842
14,41774,2014,3.633318781899142,2015,2.715291551682419
15,41957,2014,4.06235252733802,2015,4.060884847379076
16,42137,2015,3.057638025163611,2016,2.42782168586293
17,42321,2015,3.05427116350534,2016,2.270469368501771
18,42503,2016,2.319057151538662,2017,2.508920621023392
19,42689,2016,2.305809238174006,2017,2.467161166346266
20,42867,2017,2.536029549059826,2018,2.546671316138061
21,43053,2017,2.509111342826809,2018,2.480855794925163
22,43145,2018,3.025024236774643,2019,2.69389938681992
23,43235,2018,3.120740332206995,2019,2.775533179497169
24,43326,2018,3.279355759764568,2019,3.107596903291299
25,43418,2018,3.296731496509198,2019,3.221757900820066
26,43510,2019,2.945303709067959,2020,2.591074440292807
27,43600,2019,2.891533899000343,2020,2.545843589346886
28,43691,2019,2.827707622797226,2020,2.413544192054795
29,43783,2019,2.861315725866587,2020,2.631992339577627
30,43875,2020,2.552688975800033,2021,2.616345720823721
31,43966,2020,2.618329006605924,2021,2.671430903007876
32,44068,2020,1.790319754067715,2021,1.691013991470625
33,44159,2020,1.790319754067715,2021,2.153309886824961
34,44251,2021,2.130407351599706,2022,2.785334366326175
35,44341,2021,2.137626121054947,2022,2.891950990452763
36,44432,2021,2.339531676162721,2022,3.941556826710224
37,44525,2021,2.339531676162721,2022,4.667362054855917
38,44617,2022,5.037171918133976,2023,3.641364543513781
39,44706,2022,4.951039758187648,2023,3.481452844954491
40,44798,2022,4.834496776263886,2023,2.845322256798233
41,44890,2022,4.834496776263886,2023,3.305715257492858
42,44981,2023,3.153537734543965,2024,2.935215611250452
43,45071,2023,2.838865660558509,2024,2.377254777217375
44,45163,2023,2.798216547494237,2024,2.138412043368865
45,45254,2023,2.798216547494237,2024,1.757655717321982
46,45345,2024,1.831762447564067,2025,2.69124964061378
47,45436,2024,1.625773169906108,2025,2.42082970885531
48,45534,2024,1.530879676868468,2025,2.01742511619909
49,45618,2024,1.530879676868468,2025,2.159361127638926
50,45713,2025,2.104676416355189,2026,2.754798876280251
51,45800,2025,2.030491763452114,2026,2.559374235215039
52,45891,2025,2.060859685319461,2026,2.733459627814305
"@

$rows = $data -split "`n" | Where-Object { $_.Trim().Length -gt 0 }

foreach ($line in $rows) {
    $parts = $line.Trim() -split ","
    $r = [int]$parts[0]
    $aVal = $parts[1]
    $bVal = $parts[2]
    $cVal = $parts[3]
    $dVal = $parts[4]
    $eVal = $parts[5]

    $ws.Cells.Item($r, 1).Value = [double]$aVal
    $ws.Cells.Item($r, 2).Value = [double]$bVal

    if ([string]::IsNullOrEmpty($cVal)) {
        $ws.Cells.Item($r, 3).ClearContents()
    } else {
        $ws.Cells.Item($r, 3).Value = [double]$cVal
    }

    $ws.Cells.Item($r, 4).Value = [double]$dVal

    if ([string]::IsNullOrEmpty($eVal)) {
        $ws.Cells.Item($r, 5).ClearContents()
    } else {
        $ws.Cells.Item($r, 5).Value = [double]$eVal
    }
}
